# Migrated to no-subcategory model
# The "Sub-Category" column (C) values are removed (content cleared, formatting kept)
# for all data rows (2-55). The previously selected cell C10 is replaced by D24,
# and several row heights shrink slightly now that column C no longer wraps text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Sub-Category" values in column C for all data rows, keeping styles intact.
$ws.Range("C2:C55").ClearContents()

# Update row heights that changed as a consequence of the shorter column C content.
$rowsAt12_8 = @(6, 7, 9)
foreach ($r in $rowsAt12_8) {
    $ws.Rows.Item($r).RowHeight = 12.8
}

$rowsAt13_8 = 12..55
foreach ($r in $rowsAt13_8) {
    $ws.Rows.Item($r).RowHeight = 13.8
}

# Move the active selection from C10 to D24.
$ws.Range("D24").Select() | Out-Null
